$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Sheet "Trilhas": update track names/URLs and swap labels
# ---------------------------------------------------------------------------
$wsTrilhas = $wb.Worksheets.Item("Trilhas")
$wsTrilhas.Range("B2").Value2 = "Associate Data Scientist in Python"
$wsTrilhas.Range("C2").Value2 = "https://app.datacamp.com/learn/career-tracks/associate-data-scientist-in-python"
$wsTrilhas.Range("B3").Value2 = "Associate Data Engineer in SQL"
$wsTrilhas.Range("C3").Value2 = "https://app.datacamp.com/learn/career-tracks/associate-data-engineer-in-sql"
$wsTrilhas.Range("B4").Value2 = "Capacitação 2025 - Básico"

# ---------------------------------------------------------------------------
# 2) Sheet "Cursos": add a "duracao" column (C) with header + values
# ---------------------------------------------------------------------------
$wsCursos = $wb.Worksheets.Item("Cursos")
$headerSrc = $wsCursos.Range("B1")
$headerDst = $wsCursos.Range("C1")
$headerSrc.Copy($headerDst)
$headerDst.Value2 = "duracao"

$wsCursos.Range("C2").Value2 = 4
$wsCursos.Range("C3").Value2 = 4
$wsCursos.Range("C4").Value2 = 4

# ---------------------------------------------------------------------------
# 3) Sheet "Trilhas_tem_Cursos": swap id_trilha values in column A
# ---------------------------------------------------------------------------
$wsTTC = $wb.Worksheets.Item("Trilhas_tem_Cursos")
$wsTTC.Range("A2").Value2 = 0
$wsTTC.Range("A3").Value2 = 0
$wsTTC.Range("A4").Value2 = 1

# ---------------------------------------------------------------------------
# 4) New sheet "membro_feadev_faz_trilhas"
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsMFT = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$wsMFT.Name = "membro_feadev_faz_trilhas"

$headers1 = @("id_membro", "id_trilha", "data_inicio", "data_fim", "finalizado")
for ($i = 0; $i -lt $headers1.Length; $i++) {
    $col = $i + 1
    $src = $wsTrilhas.Range("A1")
    $dst = $wsMFT.Cells.Item(1, $col)
    $src.Copy($dst)
    $dst.Value2 = $headers1[$i]
}

$wsMFT.Cells.Item(2, 1).Value2 = 1
$wsMFT.Cells.Item(2, 2).Value2 = 2
$wsMFT.Cells.Item(2, 3).Value2 = "15/06/2025"
$wsMFT.Cells.Item(2, 4).Value2 = "20/06/2025"
$wsMFT.Cells.Item(2, 5).Value2 = $true

# ---------------------------------------------------------------------------
# 5) New sheet "membro_feadev_faz_cursos"
# ---------------------------------------------------------------------------
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsMFC = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet2)
$wsMFC.Name = "membro_feadev_faz_cursos"

$headers2 = @("id_membro", "id_curso", "data_inicio", "data_fim", "finalizado")
for ($i = 0; $i -lt $headers2.Length; $i++) {
    $col = $i + 1
    $src = $wsTrilhas.Range("A1")
    $dst = $wsMFC.Cells.Item(1, $col)
    $src.Copy($dst)
    $dst.Value2 = $headers2[$i]
}

$wsMFC.Cells.Item(2, 1).Value2 = 1
$wsMFC.Cells.Item(2, 2).Value2 = 0
$wsMFC.Cells.Item(2, 3).Value2 = "15/06/2025"
$wsMFC.Cells.Item(2, 4).Value2 = "15/06/2025"
$wsMFC.Cells.Item(2, 5).Value2 = $true

$wsMFC.Cells.Item(3, 1).Value2 = 1
$wsMFC.Cells.Item(3, 2).Value2 = 1
$wsMFC.Cells.Item(3, 3).Value2 = "17/06/2025"
$wsMFC.Cells.Item(3, 4).Value2 = "21/06/2025"
$wsMFC.Cells.Item(3, 5).Value2 = $true

$wsMFC.Cells.Item(4, 1).Value2 = 1
$wsMFC.Cells.Item(4, 2).Value2 = 2
$wsMFC.Cells.Item(4, 3).Value2 = "20/06/2025"
$wsMFC.Cells.Item(4, 4).Value2 = "21/06/2025"
$wsMFC.Cells.Item(4, 5).Value2 = $true

Write-Host "Done"
